$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric values (per diff) for rows 2-34
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.08186397984886642
$ws.Range("G2").Value = 0.007194244604316547
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.1813602015113354
$ws.Range("G3").Value = 0.04316546762589927
$ws.Range("H3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.003597122302158274
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.03568827385287685
$ws.Range("B5").Value = 0.1508838383838407
$ws.Range("C5").Value = 0.0239294710327456
$ws.Range("I5").Value = 0.0625
$ws.Range("J5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.05667506297229215
$ws.Range("G7").Value = 0.01798561151079137
$ws.Range("H7").Value = 0.0007434944237918215
$ws.Range("K7").Value = 0.08368200836820086
$ws.Range("B8").Value = 0.1319444444444462
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 0.6691729323308216
$ws.Range("G8").Value = 0.2805755395683454
$ws.Range("H8").Value = 0
$ws.Range("K9").Value = 0.08368200836820086
$ws.Range("E10").Value = 0
$ws.Range("H10").Value = 0.02825278810408923
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0.183175528040786
$ws.Range("E12").Value = 0
$ws.Range("H12").Value = 0.0007434944237918215
$ws.Range("J12").Value = 0.02549162418062634
$ws.Range("B13").Value = 0.01809764309764311
$ws.Range("C13").Value = 0.002518891687657431
$ws.Range("B14").Value = 0.01599326599326602
$ws.Range("C14").Value = 0
$ws.Range("I14").Value = 0.9375
$ws.Range("J14").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0.02075746540422432
$ws.Range("B16").Value = 0.04966329966329953
$ws.Range("C16").Value = 0.06801007556675058
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0.007194244604316547
$ws.Range("H16").Value = 0.01115241635687733
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0.03095411507647482
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0.01133501259445844
$ws.Range("E18").Value = 0
$ws.Range("H18").Value = 0.0007434944237918215
$ws.Range("J18").Value = 0.002184996358339403
$ws.Range("B19").Value = 0.02609427609427607
$ws.Range("C19").Value = 0.05919395465994958
$ws.Range("D19").Value = 0.006427604871447904
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0.0516372795969773
$ws.Range("B21").Value = 0.0117845117845118
$ws.Range("C21").Value = 0.03148614609571789
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0.003597122302158274
$ws.Range("H22").Value = 0.02156133828996283
$ws.Range("J22").Value = 0.03168244719592129
$ws.Range("B23").Value = 0.002946127946127945
$ws.Range("C23").Value = 0.3186397984886651
$ws.Range("G23").Value = 0.0827338129496403
$ws.Range("H23").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0.004734158776402039
$ws.Range("B32").Value = 0.01136363636363637
$ws.Range("C32").Value = 0.06926952141057929
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0.2194244604316549
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0.09686817188638054
$ws.Range("E33").Value = 0
$ws.Range("H33").Value = 0.1613382899628248
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0.1449380917698474
$ws.Range("B34").Value = 0.0006313131313131313
$ws.Range("C34").Value = 0.003778337531486146

# Remove trailing rows 36-40 (no longer present in the updated dataset)
$ws.Range("A36:K40").EntireRow.Delete()
